$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1. New "helper" columns I..O on rows 3-12: math/comparison operators and
#    cell referencing against the goal cells C15/C16.
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 12; $r++) {
    $ws.Range("I$r").Formula = "=H$r*D$r"
    $ws.Range("J$r").Formula = "=I$r=E$r"
    $ws.Range("K$r").Formula = "=E$r>D$r"
    $ws.Range("L$r").Formula = "=C$r<=`$C`$15"
    $ws.Range("M$r").Formula = "=D$r>=`$C`$16"
    $ws.Range("N$r").Formula = "=L$r*M$r"
    $ws.Range("O$r").Formula = "=N$r=1"
}

# ---------------------------------------------------------------------------
# 2. Goal values used by the new L/M columns.
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = 5
$ws.Range("C16").Value = 80000

# ---------------------------------------------------------------------------
# 3. Percent-format the "Bonus Rate (/)" column (H) header + data as percent.
#    (data cells styled first so the new-style indices line up the same way
#    Excel would have generated them)
# ---------------------------------------------------------------------------
$hData = $ws.Range("H3:H12")
$hData.NumberFormat = "0%"
$hData.Font.Name = "Calibri"
$hData.Font.Bold = $false

$h2 = $ws.Range("H2")
$h2.NumberFormat = "0%"
$h2.Font.Name = "Calibri"
$h2.Font.Bold = $true
$h2.HorizontalAlignment = -4108   # xlCenter
$h2.WrapText = $true

# ---------------------------------------------------------------------------
# 4. View state: scrolled so column B is the left-most visible column, and
#    selection resting on I9.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I9").Select() | Out-Null
